$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a full data row (columns A..R)
function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q, $rr) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rr
}

# Row 75: new weekly "Choclero" entry dated 2022-02-23 (was the "Dulce o Americano" row)
Set-Row 75 11 "Vega Monumental Concepción" "Bíobío" 44615 8 100112024 "Choclo" "Choclero" "Primera" 10000 150 160 155 "`$/unidad" "Región de O'Higgins" 155 1 "Hortaliza"

# Row 76: "Segunda" quality entry also dated 2022-02-23
Set-Row 76 11 "Vega Monumental Concepción" "Bíobío" 44615 8 100112024 "Choclo" "Choclero" "Segunda" 5000 120 120 120 "`$/unidad" "Región de O'Higgins" 120 1 "Hortaliza"

# Row 77: shifted-down "Dulce o Americano" entry dated 2021-06-02
Set-Row 77 11 "Vega Monumental Concepción" "Bíobío" 44349 8 100112024 "Choclo" "Dulce o Americano" "Primera" 50 9000 10000 9600 "`$/malla 60 unidades" "Provincia de Limarí" 160 60 "Hortaliza"

# Row 78: shifted-down "Choclero" entry dated 2022-02-03
Set-Row 78 11 "Vega Monumental Concepción" "Bíobío" 44595 8 100112024 "Choclo" "Choclero" "Primera" 20000 150 200 175 "`$/unidad" "Región de O'Higgins" 175 1 "Hortaliza"

# Row 79 (new): "Choclero" malla 50 unidades entry dated 2021-12-22
Set-Row 79 11 "Vega Monumental Concepción" "Bíobío" 44552 8 100112024 "Choclo" "Choclero" "Primera" 100 22000 24000 23000 "`$/malla 50 unidades" "Región de Arica y Parinacota" 460 50 "Hortaliza"

# Row 80 (new): "Choclero" malla 70 unidades entry dated 2021-12-22
Set-Row 80 11 "Vega Monumental Concepción" "Bíobío" 44552 8 100112024 "Choclo" "Choclero" "Primera" 100 11000 12000 11500 "`$/malla 70 unidades" "Región de Arica y Parinacota" 164 70 "Hortaliza"
